$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated capital structure database values (rows 2-10)

# Row 2
$ws.Range("D2").Value = -0.007725
$ws.Range("E2").Value = -0.1015
$ws.Range("F2").Value = -0.019
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = -5220
$ws.Range("L2").Value = -0.06410413852388555
$ws.Range("M2").Value = 6708.5834
$ws.Range("N2").Value = 0.05738609235380565
$ws.Range("O2").Value = -1.285169233716475
$ws.Range("P2").Value = 5415.6834
$ws.Range("Q2").Value = 0.04632645809417412
$ws.Range("R2").Value = -1.03748724137931
$ws.Range("S2").Value = 1292.9
$ws.Range("T2").Value = 0.1927232506344037
$ws.Range("U2").Value = 304174.4
$ws.Range("V2").Value = 2.601947262079714
$ws.Range("W2").Value = 0.01991026507688584
$ws.Range("X2").Value = 0.1336179048857245
$ws.Range("Y2").Value = -0.1137076398088387
$ws.Range("Z2").Value = 0.1275448270806661
$ws.Range("AB2").Value = 0.04812818290305465
$ws.Range("AC2").Value = -0.04812818290305465
$ws.Range("AD2").Value = 657771.8
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 657771.8
$ws.Range("AG2").Value = 353597.4
$ws.Range("AH2").Value = 0.8490945357171994
$ws.Range("AI2").Value = 0.734790809113316
$ws.Range("AJ2").Value = 0.7515353878852286
$ws.Range("AK2").Value = 0.5982954545454546

# Row 3
$ws.Range("D3").Value = -0.00813
$ws.Range("E3").Value = -0.058
$ws.Range("F3").Value = 0.161
$ws.Range("K3").Value = 106
$ws.Range("L3").Value = 0.09335094671950682
$ws.Range("U3").Value = 5401.6
$ws.Range("V3").Value = 3.909386987044945
$ws.Range("W3").Value = 0.02424242424242424
$ws.Range("X3").Value = 0.05551181641836642
$ws.Range("Y3").Value = -0.03126939217594217
$ws.Range("Z3").Value = 0.5340262427691295
$ws.Range("AB3").Value = 0.04482106437837653
$ws.Range("AC3").Value = -0.04482106437837653
$ws.Range("AD3").Value = 992.2
$ws.Range("AF3").Value = 992.2
$ws.Range("AG3").Value = -4409.400000000001
$ws.Range("AH3").Value = 0.4179620034542315
$ws.Range("AI3").Value = 0.1717411247468541
$ws.Range("AJ3").Value = 1.456353007233213
$ws.Range("AK3").Value = -11.73649188182061

# Row 4
$ws.Range("D4").Value = 0.08169999999999999
$ws.Range("E4").Value = -0.0119
$ws.Range("F4").Value = 0.00428
$ws.Range("K4").Value = 383.2
$ws.Range("L4").Value = 0.1735507246376811
$ws.Range("M4").Value = 103.3505
$ws.Range("N4").Value = 0.02125023131489668
$ws.Range("O4").Value = 0.2697038100208768
$ws.Range("P4").Value = 103.3505
$ws.Range("Q4").Value = 0.02125023131489668
$ws.Range("R4").Value = 0.2697038100208768
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0.07348175420429923
$ws.Range("X4").Value = 0.0902848682339171
$ws.Range("Y4").Value = -0.01680311402961787
$ws.Range("Z4").Value = 0.3484463522022505
$ws.Range("AB4").Value = 0.04765620961932324
$ws.Range("AC4").Value = -0.04765620961932324
$ws.Range("AD4").Value = 11000
$ws.Range("AF4").Value = 11000
$ws.Range("AG4").Value = 11000
$ws.Range("AH4").Value = 0.6934157027137769
$ws.Range("AI4").Value = 0.6649137123341493
$ws.Range("AJ4").Value = 0.6934157027137769
$ws.Range("AK4").Value = 0.6649137123341493

# Row 5
$ws.Range("D5").Value = -0.0113
$ws.Range("F5").Value = -0.0746
$ws.Range("K5").Value = -7356.7
$ws.Range("L5").Value = -0.1808956383611766
$ws.Range("M5").Value = 1958.0527
$ws.Range("N5").Value = 0.0363978054089624
$ws.Range("O5").Value = -0.2661591066646731
$ws.Range("P5").Value = 1958.0527
$ws.Range("Q5").Value = 0.0363978054089624
$ws.Range("R5").Value = -0.2661591066646731
$ws.Range("U5").Value = 173996
$ws.Range("V5").Value = 3.234372879717599
$ws.Range("W5").Value = -0.06888244282814361
$ws.Range("X5").Value = 0.1986328546522807
$ws.Range("Y5").Value = -0.2675152974804243
$ws.Range("Z5").Value = 0.1178268601503916
$ws.Range("AB5").Value = 0.04774246018771645
$ws.Range("AC5").Value = -0.04774246018771645
$ws.Range("AD5").Value = 380419.2
$ws.Range("AF5").Value = 380419.2
$ws.Range("AG5").Value = 206423.2
$ws.Range("AH5").Value = 0.8761077171199251
$ws.Range("AI5").Value = 0.7801216017770454
$ws.Range("AJ5").Value = 0.7932669046968497
$ws.Range("AK5").Value = 0.6581430146002692

# Row 6
$ws.Range("D6").Value = 0.09570000000000001
$ws.Range("E6").Value = 0.0346
$ws.Range("F6").Value = -0.033
$ws.Range("K6").Value = 1368
$ws.Range("L6").Value = 0.1527450564419781
$ws.Range("M6").Value = 490.1468
$ws.Range("N6").Value = 0.03190644447337586
$ws.Range("O6").Value = 0.3582944444444444
$ws.Range("P6").Value = 490.1468
$ws.Range("Q6").Value = 0.03190644447337586
$ws.Range("R6").Value = 0.3582944444444444
$ws.Range("U6").Value = 58723.6
$ws.Range("V6").Value = 3.822653300351517
$ws.Range("W6").Value = 0.05085293909914465
$ws.Range("X6").Value = 0.099983037474287
$ws.Range("Y6").Value = -0.04913009837514234
$ws.Range("Z6").Value = 0.2135577620083172
$ws.Range("AB6").Value = 0.04808517311852974
$ws.Range("AC6").Value = -0.04808517311852974
$ws.Range("AD6").Value = 41358.6
$ws.Range("AF6").Value = 41358.6
$ws.Range("AG6").Value = -17365
$ws.Range("AH6").Value = 0.7291636548273467
$ws.Range("AI6").Value = 0.5892554129700419
$ws.Range("AJ6").Value = 8.669495756365452
$ws.Range("AK6").Value = -1.514702162364907

# Row 7
$ws.Range("D7").Value = 0.0107
$ws.Range("F7").Value = 0.169
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = -199.6
$ws.Range("L7").Value = -0.009463797182671199
$ws.Range("M7").Value = 3455.9
$ws.Range("N7").Value = 0.1054811495824584
$ws.Range("O7").Value = -17.31412825651303
$ws.Range("P7").Value = 2163
$ws.Range("Q7").Value = 0.06601919226449186
$ws.Range("R7").Value = -10.83667334669339
$ws.Range("S7").Value = 1292.9
$ws.Range("T7").Value = 0.3741138343123355
$ws.Range("U7").Value = 13379.5
$ws.Range("V7").Value = 0.4083697563119598
$ws.Range("W7").Value = -0.003608378814694562
$ws.Range("X7").Value = 0.1433281372904136
$ws.Range("Y7").Value = -0.1469365161051082
$ws.Range("Z7").Value = 0.1220683715767042
$ws.Range("AA7").Value = 0
$ws.Range("AB7").Value = 0.04817119268757956
$ws.Range("AC7").Value = -0.04817119268757956
$ws.Range("AD7").Value = 151249.4
$ws.Range("AE7").Value = 0
$ws.Range("AF7").Value = 151249.4
$ws.Range("AG7").Value = 137869.9
$ws.Range("AH7").Value = 0.8219513228985406
$ws.Range("AI7").Value = 0.7263723930331366
$ws.Range("AJ7").Value = 0.8079903606041265
$ws.Range("AK7").Value = 0.707583211784474

# Row 8
$ws.Range("D8").Value = -0.0461
$ws.Range("E8").Value = -0.145
$ws.Range("F8").Value = -0.005
$ws.Range("K8").Value = 86.5
$ws.Range("L8").Value = 0.1324249846907532
$ws.Range("M8").Value = 23.8328
$ws.Range("N8").Value = 0.02595317434389633
$ws.Range("O8").Value = 0.2755236994219653
$ws.Range("P8").Value = 23.8328
$ws.Range("Q8").Value = 0.02595317434389633
$ws.Range("R8").Value = 0.2755236994219653
$ws.Range("T8").Value = 0
$ws.Range("U8").Value = 2430.7
$ws.Range("V8").Value = 2.646956332353261
$ws.Range("W8").Value = 0.02549140954233342
$ws.Range("X8").Value = 0.1239076724810355
$ws.Range("Y8").Value = -0.09841626293870205
$ws.Range("Z8").Value = 0.1180169111801691
$ws.Range("AB8").Value = 0.04969954694050788
$ws.Range("AC8").Value = -0.04969954694050788
$ws.Range("AD8").Value = 3447.6
$ws.Range("AF8").Value = 3447.6
$ws.Range("AG8").Value = 1016.9
$ws.Range("AH8").Value = 0.7896653610939326
$ws.Range("AI8").Value = 0.476003755453692
$ws.Range("AJ8").Value = 0.5254754030591153
$ws.Range("AK8").Value = 0.2113214604850273

# Row 9
$ws.Range("D9").Value = -0.0614
$ws.Range("E9").Value = -0.289
$ws.Range("F9").Value = -0.07000000000000001
$ws.Range("K9").Value = 172
$ws.Range("L9").Value = 0.06173725771715721
$ws.Range("M9").Value = 413.1136
$ws.Range("N9").Value = 0.07672986627043091
$ws.Range("O9").Value = 2.401823255813953
$ws.Range("P9").Value = 413.1136
$ws.Range("Q9").Value = 0.07672986627043091
$ws.Range("R9").Value = 2.401823255813953
$ws.Range("U9").Value = 18808.1
$ws.Range("V9").Value = 3.493332095096582
$ws.Range("W9").Value = 0.01179196775033936
$ws.Range("X9").Value = 0.1882167947010388
$ws.Range("Y9").Value = -0.1764248269506994
$ws.Range("Z9").Value = 0.07062032324133595
$ws.Range("AB9").Value = 0.05112227720757786
$ws.Range("AC9").Value = -0.05112227720757786
$ws.Range("AD9").Value = 35583.6
$ws.Range("AF9").Value = 35583.6
$ws.Range("AG9").Value = 16775.5
$ws.Range("AH9").Value = 0.868579072242455
$ws.Range("AI9").Value = 0.6991928067845101
$ws.Range("AJ9").Value = 0.7570342291116677
$ws.Range("AK9").Value = 0.5228569736600144

# Row 10
$ws.Range("D10").Value = -0.00732
$ws.Range("E10").Value = -0.204
$ws.Range("F10").Value = -0.159
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 220.6
$ws.Range("L10").Value = 0.05610233717352051
$ws.Range("M10").Value = 264.187
$ws.Range("N10").Value = 0.1085402629416598
$ws.Range("O10").Value = 1.197583862194016
$ws.Range("P10").Value = 264.187
$ws.Range("Q10").Value = 0.1085402629416598
$ws.Range("R10").Value = 1.197583862194016
$ws.Range("U10").Value = 31434.9
$ws.Range("V10").Value = 12.91491372226787
$ws.Range("W10").Value = 0.01557810591134744
$ws.Range("X10").Value = 0.3514240417604667
$ws.Range("Y10").Value = -0.3358459358491193
$ws.Range("Z10").Value = 0.1565027383302554
$ws.Range("AA10").Value = 0
$ws.Range("AB10").Value = 0.05199213006315113
$ws.Range("AC10").Value = -0.05199213006315113
$ws.Range("AD10").Value = 33721.2
$ws.Range("AE10").Value = 0
$ws.Range("AF10").Value = 33721.2
$ws.Range("AG10").Value = 2286.299999999996
$ws.Range("AH10").Value = 0.9326791166969066
$ws.Range("AI10").Value = 0.692826867247837
$ws.Range("AJ10").Value = 0.4843548079571209
$ws.Range("AK10").Value = 0.1326390903289433

# Cells removed entirely in the updated dataset
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("AN7").ClearContents()
$ws.Range("AP7").ClearContents()
$ws.Range("AN10").ClearContents()
$ws.Range("AP10").ClearContents()
